# Refresh the cryptocurrency price / 1h-volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    [PSCustomObject]@{ Row = 2; D = '64.203.09'; E = '  +1.39%  ' }
    [PSCustomObject]@{ Row = 3; D = '2.789.54'; E = '  +2.49%  ' }
    [PSCustomObject]@{ Row = 4; D = '0.999'; E = '  -0.08%  ' }
    [PSCustomObject]@{ Row = 5; D = '591.55'; E = '  +0.69%  ' }
    [PSCustomObject]@{ Row = 6; D = '161.01'; E = '  +6.49%  ' }
    [PSCustomObject]@{ Row = 7; D = '0.623'; E = '  +2.43%  ' }
    [PSCustomObject]@{ Row = 8; E = '  +0.11%  ' }
    [PSCustomObject]@{ Row = 9; E = '  +0.70%  ' }
    [PSCustomObject]@{ Row = 10; D = '0.115'; E = '  +1.65%  ' }
    [PSCustomObject]@{ Row = 11; D = '0.398'; E = '  +1.93%  ' }
    [PSCustomObject]@{ Row = 12; E = '  +1.06%  ' }
    [PSCustomObject]@{ Row = 13; D = '3.283.59'; E = '  +2.29%  ' }
    [PSCustomObject]@{ Row = 14; D = '27.45'; E = '  +2.32%  ' }
    [PSCustomObject]@{ Row = 15; D = '64.099.47'; E = '  +1.42%  ' }
    [PSCustomObject]@{ Row = 16; E = '  +5.70%  ' }
    [PSCustomObject]@{ Row = 17; D = '2.790.95'; E = '  +1.04%  ' }
    [PSCustomObject]@{ Row = 18; D = '12.47'; E = '  +3.87%  ' }
    [PSCustomObject]@{ Row = 19; D = '5.06'; E = '  +3.69%  ' }
    [PSCustomObject]@{ Row = 20; D = '368.06'; E = '  +1.06%  ' }
    [PSCustomObject]@{ Row = 21; D = '7.06'; E = '  +0.27%  ' }
    [PSCustomObject]@{ Row = 22; D = '0.577'; E = '  +7.58%  ' }
    [PSCustomObject]@{ Row = 23; D = '0.999'; E = '  +0.10%  ' }
    [PSCustomObject]@{ Row = 24; D = '67.38'; E = '  +2.47%  ' }
    [PSCustomObject]@{ Row = 25; E = '  +6.41%  ' }
    [PSCustomObject]@{ Row = 26; D = '8.87'; E = '  +3.28%  ' }
    [PSCustomObject]@{ Row = 27; D = '0.0₃0975'; E = '  +12.37%  ' }
    [PSCustomObject]@{ Row = 28; E = '  +0.33%  ' }
    [PSCustomObject]@{ Row = 29; E = '  +1.86%  ' }
    [PSCustomObject]@{ Row = 30; D = '7.28'; E = '  +2.08%  ' }
    [PSCustomObject]@{ Row = 31; E = '  +6.10%  ' }
    [PSCustomObject]@{ Row = 32; D = '5.21'; E = '  +9.10%  ' }
    [PSCustomObject]@{ Row = 33; D = '170.47'; E = '  -0.19%  ' }
    [PSCustomObject]@{ Row = 34; D = '20.94'; E = '  +1.70%  ' }
    [PSCustomObject]@{ Row = 35; E = '  +0.18%  ' }
    [PSCustomObject]@{ Row = 36; E = '  +4.29%  ' }
    [PSCustomObject]@{ Row = 37; E = '  +2.59%  ' }
    [PSCustomObject]@{ Row = 38; E = '  +1.96%  ' }
    [PSCustomObject]@{ Row = 39; B = 'Bittensor'; C = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D = '344.09'; E = '  -2.18%  ' }
    [PSCustomObject]@{ Row = 40; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '6.33'; E = '  +11.16%  ' }
    [PSCustomObject]@{ Row = 41; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.29'; E = '  +0.21%  ' }
    [PSCustomObject]@{ Row = 42; D = '40.30'; E = '  +2.47%  ' }
    [PSCustomObject]@{ Row = 43; D = '22.59'; E = '  +0.35%  ' }
    [PSCustomObject]@{ Row = 44; D = '0.0618'; E = '  +4.11%  ' }
    [PSCustomObject]@{ Row = 45; D = '22.56'; E = '  +2.83%  ' }
    [PSCustomObject]@{ Row = 46; D = '0.654'; E = '  +1.75%  ' }
    [PSCustomObject]@{ Row = 47; E = '  +1.68%  ' }
    [PSCustomObject]@{ Row = 48; D = '139.16'; E = '  -0.31%  ' }
    [PSCustomObject]@{ Row = 49; E = '  +2.24%  ' }
    [PSCustomObject]@{ Row = 50; D = '2.178.51'; E = '  +0.35%  ' }
    [PSCustomObject]@{ Row = 51; D = '0.998'; E = '  +0.39%  ' }
)

foreach ($u in $updates) {
    if ($u.B) {
        $ws.Cells.Item($u.Row, 2).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
    if ($u.C) {
        $ws.Cells.Item($u.Row, 3).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 3).Value = $u.C
    }
    if ($u.D) {
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.E) {
        $ws.Cells.Item($u.Row, 5).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
